# Auto-generated Excel COM-interop script applying the Titan_Profits.xlsx data refresh.
# For each affected leve row (columns H:N = market-price / profit figures), this sets the
# updated cached values per sheet. Cells that must be removed entirely (so they serialize as
# absent <c> elements, matching a row that now has no value in that column) use ClearContents().

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 11875
$ws.Range("I34").Value = 9285.714
$ws.Range("J34").Value = 30000
$ws.Range("K34").Value = 9285.714
$ws.Range("L34").Value = 30000
$ws.Range("M34").Value = -9082.714
$ws.Range("N34").Value = -30406
$ws.Range("H36").Value = 11875
$ws.Range("I36").Value = 9285.714
$ws.Range("J36").Value = 30000
$ws.Range("K36").Value = 9285.714
$ws.Range("L36").Value = 30000
$ws.Range("M36").Value = -8570.714
$ws.Range("N36").Value = -31430
$ws.Range("H51").Value = 2094.5
$ws.Range("I51").Value = 2170
$ws.Range("J51").Value = 2079.4
$ws.Range("K51").Value = 2170
$ws.Range("L51").Value = 2079.4
$ws.Range("M51").Value = -1686
$ws.Range("N51").Value = -3047.4
$ws.Range("H52").Value = 4355
$ws.Range("J52").Value = 4355
$ws.Range("L52").Value = 13065
$ws.Range("N52").Value = -13385
$ws.Range("H125").Value = 9343086
$ws.Range("I125").Value = 374.66666
$ws.Range("J125").Value = 12457323
$ws.Range("K125").Value = 3371.99994
$ws.Range("L125").Value = 112115907
$ws.Range("M125").Value = -911.9999399999997
$ws.Range("N125").Value = -112120827
$ws.Range("H127").Value = 416.375
$ws.Range("I127").Value = 289.6154
$ws.Range("J127").Value = 965.6667
$ws.Range("K127").Value = 868.8462000000001
$ws.Range("L127").Value = 2897.0001
$ws.Range("M127").Value = 4091.1538
$ws.Range("N127").Value = -12817.0001
$ws.Range("H138").Value = 372846.44
$ws.Range("I138").Value = 1545486.6
$ws.Range("J138").Value = 2539.0132
$ws.Range("K138").Value = 4636459.800000001
$ws.Range("L138").Value = 7617.0396
$ws.Range("M138").Value = -4631319.800000001
$ws.Range("N138").Value = -17897.0396

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 22853.941
$ws.Range("I32").Value = 2057.6592
$ws.Range("K32").Value = 2057.6592
$ws.Range("M32").Value = -1770.6592
$ws.Range("H45").Value = 1445.8667
$ws.Range("I45").Value = 1097.4
$ws.Range("K45").Value = 1097.4
$ws.Range("M45").Value = -720.4000000000001
$ws.Range("H61").Value = 3230.9412
$ws.Range("I61").Value = 1754.25
$ws.Range("J61").Value = 4543.5557
$ws.Range("K61").Value = 1754.25
$ws.Range("L61").Value = 4543.5557
$ws.Range("M61").Value = -1542.25
$ws.Range("N61").Value = -4967.5557
$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("M102").ClearContents()
$ws.Range("N102").ClearContents()
$ws.Range("H122").Value = 2249.75
$ws.Range("I122").Value = 2000
$ws.Range("K122").Value = 6000
$ws.Range("M122").Value = -3550
$ws.Range("H132").Value = 2375.3914
$ws.Range("I132").Value = 1949.5143
$ws.Range("K132").Value = 5848.5429
$ws.Range("M132").Value = -3318.5429
$ws.Range("H136").Value = 3230.9412
$ws.Range("I136").Value = 1754.25
$ws.Range("J136").Value = 4543.5557
$ws.Range("K136").Value = 5262.75
$ws.Range("L136").Value = 13630.6671
$ws.Range("M136").Value = -2712.75
$ws.Range("N136").Value = -18730.6671

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3058.2632
$ws.Range("I105").Value = 2757.7144
$ws.Range("J105").Value = 3899.8
$ws.Range("K105").Value = 2757.7144
$ws.Range("L105").Value = 3899.8
$ws.Range("M105").Value = -1010.7144
$ws.Range("N105").Value = -7393.8
$ws.Range("H134").Value = 1951.898
$ws.Range("I134").Value = 1570.0488
$ws.Range("K134").Value = 4710.1464
$ws.Range("M134").Value = -2175.1464

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1254.5161
$ws.Range("I31").Value = 1004.4
$ws.Range("J31").Value = 2296.6667
$ws.Range("K31").Value = 1004.4
$ws.Range("L31").Value = 2296.6667
$ws.Range("M31").Value = -709.4
$ws.Range("N31").Value = -2886.6667
$ws.Range("H34").Value = 1254.5161
$ws.Range("I34").Value = 1004.4
$ws.Range("J34").Value = 2296.6667
$ws.Range("K34").Value = 1004.4
$ws.Range("L34").Value = 2296.6667
$ws.Range("M34").Value = -802.4
$ws.Range("N34").Value = -2700.6667
$ws.Range("H99").Value = 8930679
$ws.Range("I99").Value = 31250600
$ws.Range("J99").Value = 2710
$ws.Range("K99").Value = 31250600
$ws.Range("L99").Value = 2710
$ws.Range("M99").Value = -31249102
$ws.Range("N99").Value = -5706
$ws.Range("H122").Value = 2012.4615
$ws.Range("I122").Value = 916.2
$ws.Range("J122").Value = 5666.6665
$ws.Range("K122").Value = 2748.6
$ws.Range("L122").Value = 16999.9995
$ws.Range("M122").Value = -298.6000000000004
$ws.Range("N122").Value = -21899.9995
$ws.Range("H126").Value = 8930679
$ws.Range("I126").Value = 31250600
$ws.Range("J126").Value = 2710
$ws.Range("K126").Value = 93751800
$ws.Range("L126").Value = 8130
$ws.Range("M126").Value = -93749330
$ws.Range("N126").Value = -13070

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 782.875
$ws.Range("I23").Value = 2068.4
$ws.Range("J23").Value = 198.54546
$ws.Range("K23").Value = 6205.200000000001
$ws.Range("L23").Value = 595.6363799999999
$ws.Range("M23").Value = -5970.200000000001
$ws.Range("N23").Value = -1065.63638
$ws.Range("H56").Value = 5537.6924
$ws.Range("I56").Value = 5537.6924
$ws.Range("K56").Value = 5537.6924
$ws.Range("M56").Value = -5007.6924
$ws.Range("H122").Value = 719
$ws.Range("J122").Value = 878.7857
$ws.Range("L122").Value = 7909.071300000001
$ws.Range("N122").Value = -12809.0713
$ws.Range("H128").Value = 275592
$ws.Range("I128").Value = 275592
$ws.Range("K128").Value = 826776
$ws.Range("M128").Value = -821796

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H17").Value = 9583.333000000001
$ws.Range("I17").Value = 12875
$ws.Range("K17").Value = 12875
$ws.Range("M17").Value = -12707
$ws.Range("H80").Value = 2900
$ws.Range("I80").Value = 2850
$ws.Range("J80").Value = 3000
$ws.Range("K80").Value = 2850
$ws.Range("L80").Value = 3000
$ws.Range("M80").Value = -1852
$ws.Range("N80").Value = -4996
$ws.Range("H83").Value = 2900
$ws.Range("I83").Value = 2850
$ws.Range("J83").Value = 3000
$ws.Range("K83").Value = 14250
$ws.Range("L83").Value = 15000
$ws.Range("M83").Value = -9258
$ws.Range("N83").Value = -24984
$ws.Range("H122").Value = 1236145.5
$ws.Range("I122").Value = 2223582
$ws.Range("J122").Value = 1850
$ws.Range("K122").Value = 6670746
$ws.Range("L122").Value = 5550
$ws.Range("M122").Value = -6668296
$ws.Range("N122").Value = -10450
$ws.Range("H123").Value = 10322.619
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 10322.619
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 10322.619
$ws.Range("M123").ClearContents()
$ws.Range("N123").Value = -15222.619
$ws.Range("H126").Value = 2574.8
$ws.Range("I126").Value = 1905.909
$ws.Range("J126").Value = 3100.3572
$ws.Range("K126").Value = 5717.727000000001
$ws.Range("L126").Value = 9301.071599999999
$ws.Range("M126").Value = -3247.727000000001
$ws.Range("N126").Value = -14241.0716
$ws.Range("H132").Value = 2285.2837
$ws.Range("I132").Value = 1889.5122
$ws.Range("J132").Value = 2909.3845
$ws.Range("K132").Value = 5668.536599999999
$ws.Range("L132").Value = 8728.1535
$ws.Range("M132").Value = -3138.536599999999
$ws.Range("N132").Value = -13788.1535

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 3650
$ws.Range("J20").Value = 4300
$ws.Range("L20").Value = 4300
$ws.Range("N20").Value = -4752
$ws.Range("H68").Value = 1691.826
$ws.Range("I68").Value = 1616.421
$ws.Range("J68").Value = 2050
$ws.Range("K68").Value = 1616.421
$ws.Range("L68").Value = 2050
$ws.Range("M68").Value = -867.421
$ws.Range("N68").Value = -3548
$ws.Range("H71").Value = 1691.826
$ws.Range("I71").Value = 1616.421
$ws.Range("J71").Value = 2050
$ws.Range("K71").Value = 8082.105
$ws.Range("L71").Value = 10250
$ws.Range("M71").Value = -4338.105
$ws.Range("N71").Value = -17738
$ws.Range("H122").Value = 3061.4092
$ws.Range("J122").Value = 3615.3845
$ws.Range("L122").Value = 10846.1535
$ws.Range("N122").Value = -15746.1535
$ws.Range("H132").Value = 3298.7302
$ws.Range("I132").Value = 2643.244
$ws.Range("J132").Value = 4520.3184
$ws.Range("K132").Value = 7929.732
$ws.Range("L132").Value = 13560.9552
$ws.Range("M132").Value = -5399.732
$ws.Range("N132").Value = -18620.9552
$ws.Range("H136").Value = 3697.8704
$ws.Range("I136").Value = 2598.4167
$ws.Range("J136").Value = 5896.778
$ws.Range("K136").Value = 7795.250100000001
$ws.Range("L136").Value = 17690.334
$ws.Range("M136").Value = -5245.250100000001
$ws.Range("N136").Value = -22790.334

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("I13").Value = 1100
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 1100
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = -960
$ws.Range("N13").ClearContents()
$ws.Range("H122").Value = 37785.82
$ws.Range("I122").Value = 57055.777
$ws.Range("J122").Value = 3099.9
$ws.Range("K122").Value = 171167.331
$ws.Range("L122").Value = 9299.700000000001
$ws.Range("M122").Value = -168717.331
$ws.Range("N122").Value = -14199.7
